# cv & workshop update
# Insert 4 new workshop rows (two "Jun 2023" Data Science part-2 sessions,
# two "May 2023" Data Science part-1 sessions) into the workshops table,
# pushing all later rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 2 = "Jul 2023" (Statistics with R), existing row 3 = "May 2023"
# (Statistics with R). Two new rows go right after row 2, two more go right
# after what is now row 5 (the old row 3).
$ws.Rows("3:4").Insert()
$ws.Rows("6:7").Insert()

# New row 3: Jun 2023 - Data science for experimental life sciences with R (part 2) - Eng
$ws.Range("A3").Value = "Jun 2023"
$ws.Range("B3").Value = "Data science for experimental life sciences with R (part 2)"
$ws.Range("C3").Value = "Data science for experimental life sciences with R (part 2)"
$ws.Range("D3").Value = "Eng"
$ws.Range("E3").Value = "Thünen-Institut, Braunschweig via zoom "
$ws.Range("F3").Value = "20h"
$ws.Range("G3").Value = "https://biometrie-bmel.de/Kursinhalt?Kurs=6113"

# New row 4: Jun 2023 - Data Science in den experimentellen Naturwissenschaften mit R (Teil 2) - Ger
$ws.Range("A4").Value = "Jun 2023"
$ws.Range("B4").Value = "Data Science in den experimentellen Naturwissenschaften mit R (Teil 2)"
$ws.Range("C4").Value = "Data science for experimental life sciences with R (part 2)"
$ws.Range("D4").Value = "Ger"
$ws.Range("E4").Value = "Thünen-Institut, Braunschweig via zoom "
$ws.Range("F4").Value = "20h"
$ws.Range("G4").Value = "https://biometrie-bmel.de/Kursinhalt?Kurs=6112"

# New row 6: May 2023 - Data science for experimental life sciences with R (part 1) - Eng
$ws.Range("A6").Value = "May 2023"
$ws.Range("B6").Value = "Data science for experimental life sciences with R (part 1)"
$ws.Range("C6").Value = "Data science for experimental life sciences with R (part 1)"
$ws.Range("D6").Value = "Eng"
$ws.Range("E6").Value = "Thünen-Institut, Braunschweig via zoom "
$ws.Range("F6").Value = "20h"
$ws.Range("G6").Value = "https://biometrie-bmel.de/Kursinhalt?Kurs=6111"

# New row 7: May 2023 - Data Science in den experimentellen Naturwissenschaften mit R (Teil 1) - Ger
$ws.Range("A7").Value = "May 2023"
$ws.Range("B7").Value = "Data Science in den experimentellen Naturwissenschaften mit R (Teil 1)"
$ws.Range("C7").Value = "Data science for experimental life sciences with R (part 1)"
$ws.Range("D7").Value = "Ger"
$ws.Range("E7").Value = "Thünen-Institut, Braunschweig via zoom "
$ws.Range("F7").Value = "20h"
$ws.Range("G7").Value = "https://biometrie-bmel.de/Kursinhalt?Kurs=6114"

# The saved workbook now shows plain A1 as the active cell (no special
# selection highlighted further down the sheet).
$ws.Range("A1").Select()
